$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.54
$ws.Range("G2").Value = 1.57
$ws.Range("H2").Value = 6.2
$ws.Range("I2").Value = 6.6
$ws.Range("J2").Value = 4.8
$ws.Range("K2").Value = 5.2
$ws.Range("L2").Value = 1.34
$ws.Range("N2").Value = 5.1
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 2.42
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 1.55
$ws.Range("S2").Value = 2.78
$ws.Range("T2").Value = 1.76
$ws.Range("U2").Value = 2.08
$ws.Range("V2").Value = 1.17
$ws.Range("W2").Value = 2.74
$ws.Range("X2").Value = 22
$ws.Range("Y2").Value = 29
$ws.Range("Z2").Value = 60
$ws.Range("AB2").Value = 10.5
$ws.Range("AC2").Value = 12
$ws.Range("AD2").Value = 24
$ws.Range("AE2").Value = 85
$ws.Range("AF2").Value = 10.5
$ws.Range("AG2").Value = 9.6
$ws.Range("AH2").Value = 23
$ws.Range("AI2").Value = 80
$ws.Range("AJ2").Value = 14
$ws.Range("AK2").Value = 15
$ws.Range("AL2").Value = 29
$ws.Range("AM2").Value = 110
$ws.Range("AN2").Value = 7.2
$ws.Range("AO2").Value = 85
$ws.Range("F3").Value = 5.7
$ws.Range("G3").Value = 6
$ws.Range("H3").Value = 1.57
$ws.Range("I3").Value = 1.61
$ws.Range("J3").Value = 4.7
$ws.Range("K3").Value = 5.1
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 5.7
$ws.Range("O3").Value = 1.18
$ws.Range("P3").Value = 2.6
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 2.52
$ws.Range("T3").Value = 1.64
$ws.Range("U3").Value = 2.32
$ws.Range("V3").Value = 2.62
$ws.Range("W3").Value = 1.2
$ws.Range("X3").Value = 25
$ws.Range("Y3").Value = 12
$ws.Range("Z3").Value = 12
$ws.Range("AA3").Value = 16.5
$ws.Range("AB3").Value = 28
$ws.Range("AC3").Value = 11.5
$ws.Range("AE3").Value = 14.5
$ws.Range("AF3").Value = 55
$ws.Range("AG3").Value = 22
$ws.Range("AH3").Value = 19
$ws.Range("AI3").Value = 26
$ws.Range("AJ3").Value = 160
$ws.Range("AK3").Value = 65
$ws.Range("AL3").Value = 60
$ws.Range("AM3").Value = 580
$ws.Range("AN3").Value = 55
$ws.Range("AO3").Value = 6.4
$ws.Range("F4").Value = 1.42
$ws.Range("G4").Value = 1.47
$ws.Range("H4").Value = 7.2
$ws.Range("I4").Value = 9.199999999999999
$ws.Range("J4").Value = 4.9
$ws.Range("K4").Value = 5.9
$ws.Range("L4").Value = 1.31
$ws.Range("N4").Value = 5.3
$ws.Range("O4").Value = 1.21
$ws.Range("P4").Value = 2.44
$ws.Range("Q4").Value = 1.61
$ws.Range("R4").Value = 1.57
$ws.Range("S4").Value = 2.58
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1.12
$ws.Range("W4").Value = 3.1
$ws.Range("AB4").Value = 75
$ws.Range("AC4").Value = 42
$ws.Range("AF4").Value = 46
$ws.Range("AG4").Value = 95
$ws.Range("AJ4").Value = 180
$ws.Range("AK4").Value = 160
$ws.Range("F5").Value = 2.24
$ws.Range("G5").Value = 2.4
$ws.Range("H5").Value = 3.55
$ws.Range("I5").Value = 3.85
$ws.Range("J5").Value = 3.25
$ws.Range("K5").Value = 3.5
$ws.Range("L5").Value = 1.47
$ws.Range("M5").Value = 1.09
$ws.Range("N5").Value = 3.3
$ws.Range("O5").Value = 1.41
$ws.Range("P5").Value = 1.76
$ws.Range("Q5").Value = 2.16
$ws.Range("R5").Value = 1.3
$ws.Range("S5").Value = 3.95
$ws.Range("T5").Value = 1.91
$ws.Range("U5").Value = 1.99
$ws.Range("V5").Value = 1.35
$ws.Range("W5").Value = 1.72
$ws.Range("X5").Value = 12
$ws.Range("Y5").Value = 13
$ws.Range("Z5").Value = 30
$ws.Range("AA5").Value = 190
$ws.Range("AB5").Value = 9.199999999999999
$ws.Range("AD5").Value = 16
$ws.Range("AE5").Value = 65
$ws.Range("AF5").Value = 14.5
$ws.Range("AG5").Value = 12
$ws.Range("AI5").Value = 110
$ws.Range("AJ5").Value = 34
$ws.Range("AK5").Value = 34
$ws.Range("AL5").Value = 60
$ws.Range("AN5").Value = 24
$ws.Range("F6").Value = 1.8
$ws.Range("H6").Value = 5.1
$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 3.55
$ws.Range("K6").Value = 4
$ws.Range("L6").Value = 1.46
$ws.Range("N6").Value = 3.35
$ws.Range("O6").Value = 1.37
$ws.Range("R6").Value = 1.3
$ws.Range("S6").Value = 3.9
$ws.Range("T6").Value = 1.93
$ws.Range("U6").Value = 1.87
$ws.Range("V6").Value = 1.21
$ws.Range("W6").Value = 2.12
$ws.Range("X6").Value = 13
$ws.Range("AG6").Value = 10.5
$ws.Range("AL6").Value = 80
$ws.Range("AN6").Value = 21
